$wb = $excel.ActiveWorkbook
$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")

# Row 7 on "Active" is the task being completed:
#   "refactor Edit Palette Panel (and everything under it) to WithoutHaste.Windows.GUI"
# Capture its values before removing it from the Active sheet.
$taskId       = $wsActive.Cells.Item(7, 1).Value2
$taskTitle    = $wsActive.Cells.Item(7, 2).Value2
$taskCategory = $wsActive.Cells.Item(7, 4).Value2
$taskCreated  = $wsActive.Cells.Item(7, 5).Text
$doneDate     = $taskCreated

# Remove the completed task from the Active sheet; remaining rows shift up.
$wsActive.Rows.Item(7).Delete()

# Insert a new row at the top of the Inactive sheet's data and fill it in
# with the task's info, now marked as Done with a completion date.
$wsInactive.Rows.Item(2).Insert()

$wsInactive.Cells.Item(2, 1).Value = $taskId
$wsInactive.Cells.Item(2, 2).Value = $taskTitle
$wsInactive.Cells.Item(2, 3).Value = "Done"
$wsInactive.Cells.Item(2, 4).Value = $taskCategory
# Leading apostrophe forces these to stay plain text instead of being
# auto-converted into Excel date serials.
$wsInactive.Cells.Item(2, 5).Value = "'" + $taskCreated
$wsInactive.Cells.Item(2, 6).Value = "'" + $doneDate

# The new row inherits the header row's bold formatting on insert;
# strip that back to match the rest of the data rows.
$wsInactive.Range("A2:F2").ClearFormats()
